# Replacing the Kahraman model with a new one
# - Shift every timestamp in column A (rows 2-97) forward by 4 days
# - Replace the wind-production values in column B for rows 2-66 with the
#   new model's output (rows 67-97 remain 0, unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New production values (rows 2-66) coming from the new model
$newValues = @(501,510,500,501,505,506,535,570,590,604,611,611,650,661,662,683,710,752,799,830,842,864,872,877,909,915,908,907,878,823,763,698,663,633,621,586,547,534,492,465,452,445,474,475,495,538,558,574,605,621,627,636,650,689,699,690,672,677,667,633,532,529,535,525,479)

for ($r = 2; $r -le 97; $r++) {
    # Shift the date/time serial forward by 4 days, keeping the fractional time-of-day
    $oldDate = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = $oldDate + 4

    $idx = $r - 2
    if ($idx -lt $newValues.Length) {
        $ws.Cells.Item($r, 2).Value = $newValues[$idx]
    }
}
